# Add a new announcement row (CCLI reporting) to Table1 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the Excel table ("Table1") by one row; this extends the table's
# ref/autoFilter range from A1:C5 to A1:C6.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

# Fill in the values for the new row (row 6).
$ws.Range("A6").Value = 45931
$ws.Range("B6").Value = "CCLI reporting completed as of 10/1/25 pending youth's songs and any special songs."
$ws.Range("C6").Value = "Daim ntawv ceeb toom CCLI twb ua tiav rau hnub tim 10/1/25, tshuav cov nkauj uas cov hluas hu thiab tej nkauj tshwj xeeb xwb."

# Match the formatting of the previous data row for the DATE and
# ANNOUNCEMENT columns (date format / wrap text with table borders).
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)

# The LUS TSHAJ TAWM cell on the new row uses plain wrap-text formatting
# (no fill/border) rather than the bordered table style.
$ws.Range("C6").Style = "Normal"
$ws.Range("C6").Interior.Pattern = -4142
$ws.Range("C6").Borders.LineStyle = 0
$ws.Range("C6").WrapText = $true

# Resize the row to fit the wrapped text.
$ws.Rows(6).RowHeight = 45

# Leave the selection where the editor ended up.
$ws.Range("C17").Select() | Out-Null
